# Log_of_all_Blogs.xlsx - add a new entry (Post44) to the blog log table on Sheet1.
#
# The table "Table2" currently spans B10:F53 (header row 10, data rows 11-53).
# We append one new data row (row 54) with the same look & feel (number format,
# "no border" style flags) as the last existing data row (row 53), then fill in
# the new post's data and grow the table / autofilter to include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Duplicate the formatting of the last data row (53) into the new row (54) by
# copying the whole row and inserting it (shifting nothing further down, since
# row 54 is currently blank) - this keeps cell styles identical to row 53.
$ws.Rows.Item(53).Copy()
$ws.Rows.Item(54).Insert(-4121, 0)
$excel.CutCopyMode = 0

# Fill in the new post's values (order matches how the new strings were
# originally appended to the shared-string table: Hashnode link, title, then
# the Dev.to link).
$ws.Range("B54").Value = 44
$ws.Range("E54").Value = "https://programmingport.hashnode.dev/logical-and-operator-or-shell-scripting"
$ws.Range("C54").Value = "Logical 'AND' Operator | Shell Scripting "
$ws.Range("D54").Value = $ws.Range("D53").Value()
$ws.Range("F54").Value = "https://dev.to/rahulmishra05/logical-and-operator-shell-scripting-9cg"

# Grow the table (and its autofilter) to include the new row.
$lo.Resize($ws.Range("B10:F54"))

# Match the author's selection/active cell after adding the row.
[void]$ws.Range("F54").Select()
